$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 - copy formatting (bold, border, center/top align)
# from the existing header cell H1 so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 9
$ws.Range("J6").Value = 9

$ws.Range("I7").Value = 10
$ws.Range("J7").Value = 10

$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 5
